$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new worksheet "Test CU 5 Contacter l'assoc" at the end of the
#    workbook (after the current last sheet, "Test CU 4 ...").
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "Test CU 5 Contacter l'assoc"

# ---------------------------------------------------------------------------
# 2. Fill in the new sheet's data. Column order of first-write matters: it
#    controls the order new entries land in the shared-string table, so the
#    column A narrative cells are written before the "Ok"/"Ko"/observation
#    cells that reuse already-existing shared strings.
# ---------------------------------------------------------------------------
$ws5.Range("A1").Value = "Nom du test"
$ws5.Range("B1").Value = "Résultat du test"
$ws5.Range("C1").Value = "Observations"

$ws5.Range("A2").Value = "1. L'internaute non membre ou le membre de l'association accède au formulaire de contact"
$ws5.Range("A3").Value = "2. L'internaute non membre ou le membre de l'association saisit son nom, son email et le corps du message"
$ws5.Range("A4").Value = "3. L'internaute non membre ou le membre de l'association valide ses choix"
$ws5.Range("A6").Value = "5. L'internaute non membre ou le membre de l'association reçoit la validation de son message par courriel"
$ws5.Range("A7").Value = "6. L'administrateur du site reçoit une notification qu'il a bien reçu son message"
$ws5.Range("A8").Value = "7. L'administrateur répond à l'internaute non membre ou au membre de l'assocation par courriel et lui envoie"
$ws5.Range("A9").Value = "8. L'internaute non membre ou le membre de l'association reçoit un courriel avec la réponse de l'administrateur du site"
$ws5.Range("A5").Value = "4. Le site enregistre le message de l'internaute non membre ou du membre dans la base de données"

$ws5.Range("B2").Value = "Ok"
$ws5.Range("B3").Value = "Ok"
$ws5.Range("B4").Value = "Ok"
$ws5.Range("B5").Value = "Ok"
$ws5.Range("B6").Value = "Ok"
$ws5.Range("B7").Value = "Ko"
$ws5.Range("B8").Value = "Ok"
$ws5.Range("B9").Value = "Ok"

$ws5.Range("C7").Value = "Fonctionnalité à développer"

# Header row style (bold, same look as the header on the other sheets).
$ws5.Range("A1:C1").Font.Bold = $true

# Column A narrative cells (blue, like the other sheets' step descriptions).
$ws5.Range("A2").Font.Color = 0xC07000
$ws5.Range("A3").Font.Color = 0xC07000
$ws5.Range("A4").Font.Color = 0xC07000
$ws5.Range("A5").Font.Color = 0xC07000
$ws5.Range("A6").Font.Color = 0xC07000
$ws5.Range("A7").Font.Color = 0xC07000
$ws5.Range("A8").Font.Color = 0xC07000
$ws5.Range("A9").Font.Color = 0xC07000
$ws5.Range("A2:A8").VerticalAlignment = -4108

# Ok/Ko result cells (green "Ok", red "Ko") - B7 stays default/unstyled.
$ws5.Range("B2").Font.Color = 0x50B000
$ws5.Range("B3").Font.Color = 0x50B000
$ws5.Range("B4").Font.Color = 0x50B000
$ws5.Range("B5").Font.Color = 0x50B000
$ws5.Range("B6").Font.Color = 0x50B000
$ws5.Range("B8").Font.Color = 0x50B000
$ws5.Range("B9").Font.Color = 0x50B000

# Column widths (best-fit like the other sheets).
$ws5.Columns.Item(1).ColumnWidth = 95.25
$ws5.Columns.Item(2).ColumnWidth = 13.7
$ws5.Columns.Item(3).ColumnWidth = 22.4

$ps5 = $ws5.PageSetup
$ps5.PaperSize = 9
$ps5.Orientation = 1

# ---------------------------------------------------------------------------
# 3. Style the "Ok"/"Ko" result cells on sheet 4 ("Test CU 4 ...") the same
#    way (this sheet previously had no colour coding on column B).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Font.Color = 0x50B000
$ws4.Range("B3").Font.Color = 0x50B000
$ws4.Range("B4").Font.Color = 0x50B000
$ws4.Range("B5").Font.Color = 0x50B000
$ws4.Range("B6").Font.Color = 0x0000FF
$ws4.Range("B7").Font.Color = 0x50B000
$ws4.Range("B8").Font.Color = 0x0000FF

$ps4 = $ws4.PageSetup
$ps4.PaperSize = 9
$ps4.Orientation = 1

# ---------------------------------------------------------------------------
# 4. Move the active-sheet / selection cursor: sheet 1 becomes the selected
#    tab (cursor parked at A29), sheet 4 loses the "active" flag (cursor
#    parked at A12), sheet 5's cursor is parked at C22 while not active.
# ---------------------------------------------------------------------------
$ws5.Range("C22").Select()

$ws4.Range("A12").Select()

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A29").Select()
